$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp
$ws.Range("A1").Value = "Datos actualizados a 13 de Agosto de 2020 a las 21:25"

# Row 4: Estados Unidos
$ws.Range("B4").Value = 5392057
$ws.Range("C4").Value = 31755
$ws.Range("D4").Value = 2819134
$ws.Range("E4").Value = 2403118
$ws.Range("G4").Value = 674
$ws.Range("H4").Value = 169805

# Row 6: India
$ws.Range("B6").Value = 2459613
$ws.Range("C6").Value = 64142
$ws.Range("D6").Value = 1750636
$ws.Range("E6").Value = 660833
$ws.Range("G6").Value = 1006
$ws.Range("H6").Value = 48144

# Row 15
$ws.Range("G15").Value = 18
$ws.Range("H15").Value = 41347

# Row 27
$ws.Range("B27").Value = 121028
$ws.Range("C27").Value = 184
$ws.Range("D27").Value = 107419
$ws.Range("E27").Value = 4597
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = 9012

# Row 36
$ws.Range("G36").Value = 6
$ws.Range("H36").Value = 5776

# Row 49
$ws.Range("D49").Value = 39177
$ws.Range("E49").Value = 12601

# Row 77
$ws.Range("E77").Value = 6199
$ws.Range("G77").Value = 1
$ws.Range("H77").Value = 106

# Row 104: Maldivas - refreshed stats
$ws.Range("B104").Value = 5494
$ws.Range("C104").Value = 128
$ws.Range("D104").Value = 2920
$ws.Range("E104").Value = 2553

# Malaui's case count grew enough to overtake Zimbabue, Guinea Ecuatorial and
# Hungria in the ranking, so it now sits at row 106 and those three countries
# shift down one row (their own totals are unchanged).
$ws.Range("A106").Value = "Malaui"
$ws.Range("B106").Value = 4912
$ws.Range("C106").Value = 160
$ws.Range("D106").Value = 2550
$ws.Range("E106").Value = 2209
$ws.Range("G106").Value = 1
$ws.Range("H106").Value = 153

$ws.Range("A107").Value = "Zimbabue"
$ws.Range("B107").Value = 4893
$ws.Range("C107").Value = 0
$ws.Range("D107").Value = 1620
$ws.Range("E107").Value = 3151
$ws.Range("G107").Value = 0
$ws.Range("H107").Value = 122

$ws.Range("A108").Value = "Guinea Ecuatorial"
$ws.Range("B108").Value = 4821
$ws.Range("C108").Value = 0
$ws.Range("D108").Value = 2182
$ws.Range("E108").Value = 2556
$ws.Range("G108").Value = 0
$ws.Range("H108").Value = 83

$ws.Range("A109").Value = "Hungria"
$ws.Range("B109").Value = 4813
$ws.Range("C109").Value = 45
$ws.Range("D109").Value = 3561
$ws.Range("E109").Value = 645
$ws.Range("G109").Value = 2
$ws.Range("H109").Value = 607

# Row 141
$ws.Range("B141").Value = 1432
$ws.Range("C141").Value = 105
$ws.Range("E141").Value = 982
$ws.Range("G141").Value = 2
$ws.Range("H141").Value = 55
